$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pattern for the "Float number" row (D7): no more default pattern.
$ws.Range("D7").ClearContents()

# Add three new rows describing Range / NumberRange generators.
$ws.Range("A15").Value = "Gender selection from range"
$ws.Range("B15").Value = "Range"
$ws.Range("C15").Value = "y"
$ws.Range("D15").Value = "M,F"

$ws.Range("A16").Value = "Locale selection from range"
$ws.Range("B16").Value = "Range"
$ws.Range("C16").Value = "y"
$ws.Range("D16").Value = "US, UK, EU"

$ws.Range("A17").Value = "Number range"
$ws.Range("B17").Value = "NumberRange"
$ws.Range("C17").Value = "y"
$ws.Range("D17").Value = "-5:10"

# Match the look of the existing table rows (11:14): Arial 10, thin left/right
# borders, text number format on column C/D.
foreach ($addr in @("A15","B15","C15","D15","A16","B16","C16","D16","A17","B17","C17","D17")) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Arial"
  $c.Font.Size = 10
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(10).LineStyle = 1
}

foreach ($addr in @("C15","D15","C16","D16","C17","D17")) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D7").Select()
